# Coalesce operator order test: add two more computed columns (ThreeItem,
# ThreeBack) to the query-table results so the effect of operand order on
# `??` can be compared, mirroring the Power Query changes:
#   #"Added Custom2" = Table.AddColumn(#"Added Custom1", "ThreeItem", each [Account]??[Debit]??[Credit]),
#   #"Added Custom3" = Table.AddColumn(#"Added Custom2", "ThreeBack",  each [Math]??[Credit]??[Debit])

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grow the query result table (_tExample_) from G4:J6 to G4:L6, adding
#     the two new columns, then name + populate them. -----------------------
$lo = $ws.ListObjects.Item(2)
$lo.Resize($ws.Range("G4:L6"))

$ws.Range("K4").Value = "ThreeItem"
$ws.Range("L4").Value = "ThreeBack"

$ws.Range("K5").Value = "A"
$ws.Range("L5").Value = -2

$ws.Range("K6").Value = "B"
$ws.Range("L6").Value = 3

# Give the two new columns a sensible width, matching the other result cols.
$ws.Range("K1:L100").ColumnWidth = 12.21875

# --- The old "Math" coloured-font style on G5/G6 was dropped when the
#     table was reformatted for the extra columns. --------------------------
$ws.Range("G5:G6").ClearFormats()

# --- The hidden ExternalData_1 name (driven by the query table) now spans
#     the two extra columns too. ---------------------------------------------
$ws.Names.Item("ExternalData_1").RefersTo = "=Report!`$G`$4:`$L`$6"

# --- Leave the cursor where the author left it while checking the results. -
$ws.Range("O13").Select()
